$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "数学" column entirely (shrinks dimension to column B)
$ws.Range("C1:C5").Delete()

# Remove the old A1 header cell completely (no longer used, no leftover style)
$ws.Range("A1").Clear()

# New combined header in B1 (keeps the bold/centered/bordered header style from B1's original format)
$ws.Range("B1").Value = "姓名,语文,数学"

# Rows 2-5: column A becomes a numeric row index (0-based), column B holds
# the comma-joined "name,chinese,math" string (plain, unstyled like before)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "小明,90,92"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "小红,98,87"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "小刚,87,90"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "小丽,90,98"
